$wb = $excel.ActiveWorkbook

# The localization run moved on: cells that used to read "Ready for handoff"
# are now "In Translation". This string shows up as the per-locale status on
# the Overview sheet (columns E = zh-cn, F = de-de) and as the Status column
# (C) on each per-locale detail sheet (zh-cn, de-de).

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# The status columns are narrower now that "In Translation" is shorter than
# "Ready for handoff" - re-fit them to the new text (12.5 chars is the
# width that resolves to the correct fitted column width after Excel's
# internal character-width rounding).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
